# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1) onto the new
# header cell (H1) so the new column matches the style of the others
# (bold font, border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new header text and the data value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
